$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.561.62'
$ws.Range("E2").Value = '  -0.02%  '

$ws.Range("D3").Value = '1.874.05'
$ws.Range("E3").Value = '  -0.85%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.63'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.20%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9997'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.03%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4727'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.16%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2893'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.32%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06471'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.20%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.00'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.73%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07724'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.80%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7402'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.00%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '96.07'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.30%  '

$ws.Range("D14").Value = '1.868.04'
$ws.Range("E14").Value = '  -1.15%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.170'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.17%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '274.47'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.34%  '

$ws.Range("D17").Value = '30.625.56'
$ws.Range("E17").Value = '  +0.25%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.25'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.94%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.02%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007476'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.95%  '

$ws.Range("D21").Value = '2.116.91'
$ws.Range("E21").Value = '  -1.25%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.11%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.211'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.00%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.171'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.86%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '165.27'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.02%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.185'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.46%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.67'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.18%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.900'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.91%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.09960'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.30%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.346'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.62%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.507'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.47%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.234'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.56%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.084'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.93%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04763'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.25%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.119'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.37%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6914'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.98%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.717'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.10%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01845'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.73%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.756'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.38%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.249'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.11%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '73.38'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.46%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.970'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.11%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.000'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.02%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4158'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.03%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8336'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.85%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.42'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.46%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.331'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.41%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '35.30'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.36%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.970'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.64%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '915.09'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.71%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05660'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.02%  '

